# Updates the 1h crypto snapshot (price + volume-change columns) to the
# latest GitHub Actions scrape. Rows 40/41 also swap places (Frax now
# ranks above Aptos), so both rows' full contents are replaced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a target cell plus its new text. "Numeric-looking" Price
# values (plain decimals Excel would otherwise coerce to a Number) are
# flagged so we can force them back to text, matching the original
# inline-string cells (e.g. "0.9580" must stay "0.9580", not become 0.958).
$updates = @(
    @{ Cell = 'D2'; Value = '20.540.35'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +1.40%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.472.90'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +1.98%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.11%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '0.9575'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +5.00%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '277.46'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -0.32%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.3615'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -1.26%  '; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -1.51%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '39.65'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  +1.31%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '1.070'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +4.57%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.06658'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +1.86%  '; ForceText = $false }
    @{ Cell = 'E12'; Value = '  +0.18%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '5.522'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +2.27%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '18.16'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  +2.67%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '6.174'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +1.59%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '0.9580'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  +1.63%  '; ForceText = $false }
    @{ Cell = 'E17'; Value = '  +0.93%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '1.473.16'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +1.91%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '0.05928'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +5.12%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '68.85'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +0.53%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '5.490'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +1.52%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '14.54'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  +0.62%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '11.16'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  +2.63%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '2.264'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +0.58%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '20.545.10'; ForceText = $false }
    @{ Cell = 'E25'; Value = '  +1.31%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '143.23'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +3.86%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '2.130'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -2.03%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  +0.75%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '1.632.44'; ForceText = $false }
    @{ Cell = 'D30'; Value = '113.80'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +3.35%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '3.897'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +1.82%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '4.968'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  +2.39%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '0.08011'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +3.99%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '0.8053'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.516'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +4.49%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '1.215'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +6.50%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '0.05758'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -3.23%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '4.728'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.80%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.02057'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +2.84%  '; ForceText = $false }
    @{ Cell = 'B40'; Value = 'Frax'; ForceText = $false }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; ForceText = $false }
    @{ Cell = 'D40'; Value = '0.9585'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +2.98%  '; ForceText = $false }
    @{ Cell = 'B41'; Value = 'Aptos'; ForceText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false }
    @{ Cell = 'D41'; Value = '10.39'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +2.02%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.1873'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +1.70%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '7.424'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  +3.31%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '0.5275'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  +0.65%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '3.519'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -0.24%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '12.15'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +0.46%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '118.47'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -0.71%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '0.5205'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +1.12%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '1.817'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  +3.14%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.06471'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +2.07%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '0.9865'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -0.52%  '; ForceText = $false }
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        # Leading apostrophe forces Excel to store the value as text
        # (same as a user typing '0.9580 into the cell).
        $ws.Range($u.Cell).Value = "'" + $u.Value
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
